$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / rId1) - index 1
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 835
$ws1.Range("F5").Value = 1168
$ws1.Range("F6").Value = 1
$ws1.Range("F7").Value = 3783
$ws1.Range("F8").Value = 2513
$ws1.Range("F10").Value = 2361
$ws1.Range("F11").Value = 243
$ws1.Range("F14").Value = 1620
$ws1.Range("F15").Value = 632
$ws1.Range("F16").Value = 88
$ws1.Range("F17").Value = 294
$ws1.Range("F21").Value = 63
$ws1.Range("F22").Value = 426
$ws1.Range("F23").Value = 23
$ws1.Range("F24").Value = 85
$ws1.Range("F25").Value = 467
$ws1.Range("F26").Value = 662
$ws1.Range("F27").Value = 76
$ws1.Range("F28").Value = 68
$ws1.Range("F29").Value = 355
$ws1.Range("F31").Value = 1604
$ws1.Range("F32").Value = 813
$ws1.Range("F33").Value = 847
$ws1.Range("F34").Value = 1920
$ws1.Range("F35").Value = 212
$ws1.Range("F36").Value = 501
$ws1.Range("F37").Value = 76
$ws1.Range("F38").Value = 563
$ws1.Range("F39").Value = 1193
$ws1.Range("F41").Value = 402

# Sheet "全部类型" (sheet4 / rId4) - index 4
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 835
$ws4.Range("F5").Value = 1168
$ws4.Range("F6").Value = 1
$ws4.Range("F7").Value = 3783
$ws4.Range("F8").Value = 2513
$ws4.Range("F10").Value = 2361
$ws4.Range("F11").Value = 243
$ws4.Range("F14").Value = 1620
$ws4.Range("F15").Value = 632
$ws4.Range("F16").Value = 88
$ws4.Range("F17").Value = 294
$ws4.Range("F21").Value = 63
$ws4.Range("F22").Value = 426
$ws4.Range("F23").Value = 23
$ws4.Range("F24").Value = 85
$ws4.Range("F25").Value = 467
$ws4.Range("F26").Value = 662
$ws4.Range("F27").Value = 76
$ws4.Range("F31").Value = 68
$ws4.Range("F32").Value = 355
$ws4.Range("F34").Value = 1604
$ws4.Range("F35").Value = 813
$ws4.Range("F37").Value = 847
$ws4.Range("F38").Value = 1920
$ws4.Range("F39").Value = 212
$ws4.Range("F43").Value = 501
$ws4.Range("F44").Value = 76
$ws4.Range("F45").Value = 563
$ws4.Range("F46").Value = 1193
$ws4.Range("F48").Value = 402
